$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 829, pushing the existing 829-870 block
# (2026/12/29 .. 2027/01/05) down to 831-872.
$ws.Rows("829:830").Insert()

# New row 829: 2026/02/21 (Sat)
# Force the date-looking text to be stored as plain text (not auto-parsed
# as a date serial) by pre-formatting the cell as Text, then strip the
# number format back off so the cell matches the plain/unstyled cells
# around it.
$ws.Range("A829").NumberFormat = "@"
$ws.Range("A829").Value = "2026/02/21"
$ws.Range("A829").ClearFormats()
$ws.Range("B829").Value = "土"
$ws.Range("C829").Value = 22
$ws.Range("D829").Value = 39

# New row 830: 2026/02/22 (Sun)
$ws.Range("A830").NumberFormat = "@"
$ws.Range("A830").Value = "2026/02/22"
$ws.Range("A830").ClearFormats()
$ws.Range("B830").Value = "日"
$ws.Range("C830").Value = 2
$ws.Range("D830").Value = 42
